$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update laserKerf value (B6): 0.1 -> 0.12
$ws.Range("B6").Value = 0.12

# Widen column A slightly: 21.5703125 -> 21.625
# (the engine quantizes ColumnWidth to 1/6-character steps, so 125/6 is the
# closest representable value to the target 21.625)
$ws.Columns("A").ColumnWidth = 20.833333333

# Move the active selection from G5 to B7
$ws.Range("B7").Select()
